$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-16 (columns B-I) and add new column J formulas
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 117
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 26.63504923173846
$ws.Range("J2").Formula = "=SUM(B2:E2)*60/SUM(F2:G2)"

$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 146
$ws.Range("G3").Value = 80
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 35.15970444517611
$ws.Range("J3").Formula = "=SUM(B3:E3)*60/SUM(F3:G3)"

$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 296
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 38.70308501035772
$ws.Range("J4").Formula = "=SUM(B4:E4)*60/SUM(F4:G4)"

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 365
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 46.29175143652816
$ws.Range("J5").Formula = "=SUM(B5:E5)*60/SUM(F5:G5)"

$ws.Range("B6").Value = 40
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 54
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 692
$ws.Range("G6").Value = 303
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 46.3008862208055
$ws.Range("J6").Formula = "=SUM(B6:E6)*60/SUM(F6:G6)"

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 102
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 26.63504923173846
$ws.Range("J7").Formula = "=SUM(B7:E7)*60/SUM(F7:G7)"

$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 156
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 35.15970444517611
$ws.Range("J8").Formula = "=SUM(B8:E8)*60/SUM(F8:G8)"

$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 248
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 33.792245485229586
$ws.Range("J9").Formula = "=SUM(B9:E9)*60/SUM(F9:G9)"

$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 433
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 43.860625911465526
$ws.Range("J10").Formula = "=SUM(B10:E10)*60/SUM(F10:G10)"

$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 523
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 49.6875535888524
$ws.Range("J11").Formula = "=SUM(B11:E11)*60/SUM(F11:G11)"

$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 93
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 26.63504923173846
$ws.Range("J12").Formula = "=SUM(B12:E12)*60/SUM(F12:G12)"

$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 148
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 35.15970444517611
$ws.Range("J13").Formula = "=SUM(B13:E13)*60/SUM(F13:G13)"

$ws.Range("B14").Value = 2
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 332
$ws.Range("G14").Value = 143
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 44.43353845096168
$ws.Range("J14").Formula = "=SUM(B14:E14)*60/SUM(F14:G14)"

$ws.Range("B15").Value = 9
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 297
$ws.Range("G15").Value = 89
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 45.11794255226339
$ws.Range("J15").Formula = "=SUM(B15:E15)*60/SUM(F15:G15)"

$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 435
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 48.570820150625806
$ws.Range("J16").Formula = "=SUM(B16:E16)*60/SUM(F16:G16)"

# New header for column J
$ws.Range("J1").Value = "accidents/hour"

# Column widths (best-fit, matching AutoFit Column Width sizing from the source workbook)
$ws.Columns.Item(1).ColumnWidth = 6.0
$ws.Columns.Item(2).ColumnWidth = 10.0
$ws.Columns.Item(3).ColumnWidth = 8.166666666666666
$ws.Columns.Item(4).ColumnWidth = 13.833333333333334
$ws.Columns.Item(5).ColumnWidth = 12.0
$ws.Columns.Item(6).ColumnWidth = 19.0
$ws.Columns.Item(7).ColumnWidth = 19.833333333333332
$ws.Columns.Item(8).ColumnWidth = 13.5
$ws.Columns.Item(9).ColumnWidth = 17.666666666666668
$ws.Columns.Item(10).ColumnWidth = 12.5

# Selection moved to M11 in the saved view
[void]$ws.Range("M11").Select()
